$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: replace existing content, mark C2 as wrapped ---
$ws.Range("A2").Value = "英语Mix.png"
$ws.Range("B2").Value = "Mix"
$ws.Range("C2").Value = "Taiga Faiya Saiba Faiba Daiba Baiba JyaJya（Faibo Waipa）"
$ws.Range("D2").Value = "2八拍，语言Mix"
$ws.Range("C2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 18

# --- Row 3: new row ---
$ws.Range("A3").Value = "日语Mix.png"
$ws.Range("B3").Value = "Mix"
$ws.Range("C3").Value = "Tora Hi Jinzou Seni Ama Shindou Kasen（Tobi Jyokyo）"
$ws.Range("D3").Value = "2八拍，语言Mix"

# --- Row 5: new row (row 4 intentionally skipped) ---
$c5 = "JinzouFaiya FaiboWaipa" + [char]10 + "Taiga Taiga TatatataTaiga" + [char]10 + "ChapeApeKaraKina ChapeApeKaraKina" + [char]10 + "Myo-hontousuke (P) Waipa" + [char]10 + "Faiya Faiya ToraToraKaraKina" + [char]10 + "ChapeApeFama AmaAmaJyasupa" + [char]10 + "ToraTaiga ToraTaiga" + [char]10 + "JinzouSen'i Yetaiga"
$ws.Range("A5").Value = "可变三连.png"
$ws.Range("B5").Value = "Mix"
$ws.Range("C5").Value = $c5
$ws.Range("D5").Value = "8八拍，可变Mix"
$ws.Range("C5").WrapText = $true
$ws.Rows.Item(5).RowHeight = 141

# --- Row 6: new row ---
$c6 = "Tsukino Hikarini Terasarete" + [char]10 + "Amai Egaoni Koioshita" + [char]10 + "Hoshiga Mahoude Kagayaite" + [char]10 + "Hitomini Utsuru（Pink）iro" + [char]10 + "Sekaino Aini Tsutsumarete" + [char]10 + "Umaretekita Anatano Sobade" + [char]10 + "Zutto Isshoni Waratte Ittai" + [char]10 + "Aio Motto（Team OOO）"
$ws.Range("A6").Value = "巧克力口上.png"
$ws.Range("B6").Value = "口上"
$ws.Range("C6").Value = $c6
$ws.Range("D6").Value = "8八拍，可变Mix"
$ws.Range("C6").WrapText = $true
$ws.Rows.Item(6).RowHeight = 141

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 20.428571428571427
$ws.Columns.Item(3).ColumnWidth = 58.0
$ws.Columns.Item(4).ColumnWidth = 17.428571428571427
$ws.Columns.Item(5).ColumnWidth = 27.71428571428571

Write-Output "done"
